$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 500653.5
$ws.Range("J3").Value = 500653.5
$ws.Range("L3").Value = 500653.5
$ws.Range("N3").Value = -500881.5

$ws.Range("H12").Value = 2265.3333
$ws.Range("I12").Value = 398.5
$ws.Range("J12").Value = 5999
$ws.Range("K12").Value = 398.5
$ws.Range("L12").Value = 5999
$ws.Range("M12").Value = -228.5
$ws.Range("N12").Value = -6339

$ws.Range("H64").Value = 4832.1665
$ws.Range("J64").Value = 5098.8
$ws.Range("L64").Value = 5098.8
$ws.Range("N64").Value = -5594.8

$ws.Range("H67").Value = 4832.1665
$ws.Range("J67").Value = 5098.8
$ws.Range("L67").Value = 5098.8
$ws.Range("N67").Value = -6814.8

$ws.Range("H70").Value = 3154.8
$ws.Range("J70").Value = 3258
$ws.Range("L70").Value = 9774
$ws.Range("N70").Value = -10314

$ws.Range("H73").Value = 3154.8
$ws.Range("J73").Value = 3258
$ws.Range("L73").Value = 9774
$ws.Range("N73").Value = -11646

$ws.Range("H76").Value = 949
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

$ws.Range("H79").Value = 949
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

$ws.Range("H95").Value = 120624
$ws.Range("J95").Value = 120624
$ws.Range("L95").Value = 120624
$ws.Range("N95").Value = -126116

$ws.Range("H102").Value = 500653.5
$ws.Range("J102").Value = 500653.5
$ws.Range("L102").Value = 500653.5
$ws.Range("N102").Value = -507143.5

$ws.Range("H116").Value = 5066
$ws.Range("J116").Value = 6266
$ws.Range("L116").Value = 6266
$ws.Range("N116").Value = -13150

$ws.Range("H132").Value = 3343.7
$ws.Range("I132").Value = 3381.889
$ws.Range("K132").Value = 10145.667
$ws.Range("M132").Value = -7615.667000000001

$ws.Range("H137").Value = 4815.636
$ws.Range("I137").Value = 3357
$ws.Range("J137").Value = 7368.25
$ws.Range("K137").Value = 10071
$ws.Range("L137").Value = 22104.75
$ws.Range("M137").Value = -7521
$ws.Range("N137").Value = -27204.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2837.5
$ws.Range("J14").Value = 3250
$ws.Range("L14").Value = 3250
$ws.Range("N14").Value = -3600

$ws.Range("H19").Value = 66964.336
$ws.Range("I19").Value = 66964.336
$ws.Range("K19").Value = 66964.336
$ws.Range("M19").Value = -66735.336

$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("K25").Value = 500
$ws.Range("M25").Value = -98

$ws.Range("H43").Value = 42000
$ws.Range("I43").Value = 45000
$ws.Range("K43").Value = 45000
$ws.Range("M43").Value = -44687

$ws.Range("H45").Value = 13000
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = $null

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = $null
$ws.Range("N111").Value = 0

$ws.Range("H132").Value = 1333
$ws.Range("I132").Value = 1333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -1469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J94").Value = 9999
$ws.Range("L94").Value = 9999
$ws.Range("N94").Value = -10901

$ws.Range("H97").Value = 9981
$ws.Range("I97").Value = 9981
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 9981
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("N97").Value = -8990

$ws.Range("H99").Value = 2403
$ws.Range("J99").Value = 3673.6667
$ws.Range("L99").Value = 3673.6667
$ws.Range("N99").Value = -6669.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null

$ws.Range("H28").Value = 10663
$ws.Range("J28").Value = 10663
$ws.Range("L28").Value = 10663
$ws.Range("N28").Value = -11153

$ws.Range("H31").Value = 2425.5625
$ws.Range("J31").Value = 3225
$ws.Range("L31").Value = 3225
$ws.Range("N31").Value = -3815

$ws.Range("H34").Value = 2425.5625
$ws.Range("J34").Value = 3225
$ws.Range("L34").Value = 3225
$ws.Range("N34").Value = -3629

$ws.Range("H105").Value = 3876.25
$ws.Range("I105").Value = 2499.5
$ws.Range("K105").Value = 2499.5
$ws.Range("M105").Value = -752.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1006.5
$ws.Range("I81").Value = 1006.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3019.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -1896.5

$ws.Range("H84").Value = 1006.5
$ws.Range("I84").Value = 1006.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9058.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -3442.5

$ws.Range("H92").Value = 207.2
$ws.Range("J92").Value = 215.25
$ws.Range("L92").Value = 645.75
$ws.Range("N92").Value = -3141.75

$ws.Range("H93").Value = 13749.667
$ws.Range("I93").Value = 499
$ws.Range("J93").Value = 16399.8
$ws.Range("K93").Value = 1497
$ws.Range("L93").Value = 49199.39999999999
$ws.Range("M93").Value = 375
$ws.Range("N93").Value = -52943.39999999999

$ws.Range("H137").Value = 2710.5
$ws.Range("J137").Value = 3181.1667
$ws.Range("L137").Value = 9543.500100000001
$ws.Range("N137").Value = -19743.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 90581
$ws.Range("I62").Value = 90077
$ws.Range("K62").Value = 90077
$ws.Range("M62").Value = -89391

$ws.Range("H65").Value = 90581
$ws.Range("I65").Value = 90077
$ws.Range("K65").Value = 270231
$ws.Range("M65").Value = -266799

$ws.Range("H80").Value = 1499.1666
$ws.Range("I80").Value = 1865
$ws.Range("J80").Value = 1133.3334
$ws.Range("K80").Value = 1865
$ws.Range("L80").Value = 1133.3334
$ws.Range("M80").Value = -867
$ws.Range("N80").Value = -3129.3334

$ws.Range("H83").Value = 1499.1666
$ws.Range("I83").Value = 1865
$ws.Range("J83").Value = 1133.3334
$ws.Range("K83").Value = 9325
$ws.Range("L83").Value = 5666.666999999999
$ws.Range("M83").Value = -4333
$ws.Range("N83").Value = -15650.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4924.5
$ws.Range("I40").Value = 4566
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 4566
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -4430
$ws.Range("N40").Value = -6272

$ws.Range("H68").Value = 2051
$ws.Range("J68").Value = 2051
$ws.Range("L68").Value = 2051
$ws.Range("N68").Value = -3549

$ws.Range("H71").Value = 2051
$ws.Range("J71").Value = 2051
$ws.Range("L71").Value = 10255
$ws.Range("N71").Value = -17743

$ws.Range("H116").Value = 99680
$ws.Range("J116").Value = 99680
$ws.Range("L116").Value = 99680
$ws.Range("N116").Value = -108858

$ws.Range("H122").Value = 7473.154
$ws.Range("I122").Value = 7284.8335
$ws.Range("J122").Value = 7634.5713
$ws.Range("K122").Value = 21854.5005
$ws.Range("L122").Value = 22903.7139
$ws.Range("M122").Value = -19404.5005
$ws.Range("N122").Value = -27803.7139

$ws.Range("H132").Value = 3024.1
$ws.Range("I132").Value = 2915.6667
$ws.Range("K132").Value = 8747.000100000001
$ws.Range("M132").Value = -6217.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 911104.2
$ws.Range("I81").Value = 1561.75
$ws.Range("J81").Value = 1430842.8
$ws.Range("K81").Value = 3123.5
$ws.Range("L81").Value = 2861685.6
$ws.Range("M81").Value = -2062.5
$ws.Range("N81").Value = -2863807.6

$ws.Range("H84").Value = 911104.2
$ws.Range("I84").Value = 1561.75
$ws.Range("J84").Value = 1430842.8
$ws.Range("K84").Value = 15617.5
$ws.Range("L84").Value = 14308428
$ws.Range("M84").Value = -10313.5
$ws.Range("N84").Value = -14319036

$ws.Range("H96").Value = 1290.091
$ws.Range("I96").Value = 1313.5714
$ws.Range("J96").Value = 1249
$ws.Range("K96").Value = 1313.5714
$ws.Range("L96").Value = 1249
$ws.Range("M96").Value = 59.42859999999996
$ws.Range("N96").Value = -3995
